# Updates cryptos list values (Price / Volume(1h) columns, and a couple of
# row re-ordering / relabeling fixes) to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.743.30"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.603.32"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "211.82"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.828.89"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.603.91"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "65.07"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.10"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "143.79"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "15.36"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "3.27"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "1.287.79"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  +16.90%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "0.588"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.20"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "62.65"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "1.740.55"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "90.42"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("E51").Value = "  +0.14%  "
